# New crime data collected - weekly CompStat update (83rd Precinct)
# Moves the report forward one week: Volume 29 Number 49 -> 50,
# week of 12/5/2022-12/11/2022 -> 12/12/2022-12/18/2022, and refreshes
# every weekly/28-day/YTD/2-year figure in the crime table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number formats used by the existing numeric columns (pulled from the
# sheet's own styles so we reuse the same style slots Excel already has).
$fmtCount = "#,##0"
$fmtPct   = '#,##0.0;"-"#,##0.0'

function Set-NumCell {
    param($ws, $ref, $val)
    $ws.Range($ref).NumberFormat = $fmtCount
    $ws.Range($ref).Value2 = $val
}

function Set-PctCell {
    param($ws, $ref, $val)
    $ws.Range($ref).NumberFormat = $fmtPct
    $ws.Range($ref).Value2 = $val
}

# C14 keeps the literal text "0" for the whole script and N30 keeps the
# literal text "***.*" for the whole script, so they're safe, stable
# sources to stamp that same text + style onto other cells.
function Set-TextZero {
    param($ws, $ref)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value2 = "0"
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

function Set-TextStar {
    param($ws, $ref)
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value2 = "***.*"
    $ws.Range("N30").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

# ---- Header ----
$ws.Range("A8").Value2 = "Volume 29   Number  50"
$ws.Range("C9").Value2 = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# ---- Row 14 (Murder) ----
Set-NumCell $ws "D14" 1
Set-PctCell $ws "E14" -100
Set-NumCell $ws "G14" 1
Set-PctCell $ws "H14" -100
Set-NumCell $ws "J14" 7
Set-PctCell $ws "K14" -57.142857142857

# ---- Row 15 (Rape) ----
Set-TextZero $ws "F15"
Set-NumCell $ws "G15" 2
Set-PctCell $ws "H15" -100
Set-PctCell $ws "N15" -40.322580645161

# ---- Row 16 (Robbery) ----
Set-NumCell $ws "C16" 6
Set-NumCell $ws "D16" 6
Set-PctCell $ws "E16" 0
Set-NumCell $ws "F16" 25
Set-NumCell $ws "G16" 20
Set-PctCell $ws "H16" 25
Set-NumCell $ws "I16" 314
Set-NumCell $ws "J16" 240
Set-PctCell $ws "K16" 30.833333333333
Set-PctCell $ws "L16" 27.642276422764
Set-PctCell $ws "M16" -27.649769585253
Set-PctCell $ws "N16" -78.798109385550

# ---- Row 17 (Fel. Assault) ----
Set-NumCell $ws "C17" 3
Set-PctCell $ws "E17" -62.5
Set-NumCell $ws "F17" 23
Set-NumCell $ws "G17" 29
Set-PctCell $ws "H17" -20.689655172413
Set-NumCell $ws "I17" 362
Set-NumCell $ws "J17" 314
Set-PctCell $ws "K17" 15.286624203821
Set-PctCell $ws "L17" 19.471947194719
Set-PctCell $ws "M17" -0.549450549450
Set-PctCell $ws "N17" -55.853658536585

# ---- Row 18 (Burglary) ----
Set-NumCell $ws "C18" 6
Set-NumCell $ws "D18" 4
Set-PctCell $ws "E18" 50
Set-NumCell $ws "F18" 24
Set-NumCell $ws "G18" 20
Set-PctCell $ws "H18" 20
Set-NumCell $ws "I18" 293
Set-NumCell $ws "J18" 258
Set-PctCell $ws "K18" 13.565891472868
Set-PctCell $ws "L18" -12.275449101796
Set-PctCell $ws "M18" -35.745614035087
Set-PctCell $ws "N18" -77.769347496206

# ---- Row 19 (Gr. Larceny) ----
Set-NumCell $ws "C19" 14
Set-NumCell $ws "D19" 18
Set-PctCell $ws "E19" -22.222222222222
Set-NumCell $ws "F19" 48
Set-NumCell $ws "G19" 58
Set-PctCell $ws "H19" -17.241379310344
Set-NumCell $ws "I19" 680
Set-NumCell $ws "J19" 554
Set-PctCell $ws "K19" 22.743682310469
Set-PctCell $ws "L19" 27.579737335834
Set-PctCell $ws "M19" 130.508474576271
Set-PctCell $ws "N19" 24.087591240875

# ---- Row 20 (G.L.A.) ----
Set-NumCell $ws "C20" 5
Set-PctCell $ws "E20" 66.666666666666
Set-NumCell $ws "F20" 15
Set-NumCell $ws "G20" 13
Set-PctCell $ws "H20" 15.384615384615
Set-NumCell $ws "I20" 239
Set-NumCell $ws "J20" 169
Set-PctCell $ws "K20" 41.420118343195
Set-PctCell $ws "L20" 41.420118343195
Set-PctCell $ws "M20" 43.975903614457
Set-PctCell $ws "N20" -74.162162162162

# ---- Row 21 (TOTAL) ----
Set-NumCell $ws "C21" 34
Set-NumCell $ws "D21" 40
Set-PctCell $ws "E21" -15
Set-NumCell $ws "F21" 135
Set-NumCell $ws "G21" 143
Set-PctCell $ws "H21" -5.594405594405
Set-NumCell $ws "I21" 1928
Set-NumCell $ws "J21" 1570
Set-PctCell $ws "K21" 22.802547770700
Set-PctCell $ws "L21" 19.677219118559
Set-PctCell $ws "M21" 9.421112372304
Set-PctCell $ws "N21" -62.880246438197

# ---- Row 22 (Transit) ----
Set-NumCell $ws "D22" 1
Set-PctCell $ws "E22" -100
Set-NumCell $ws "F22" 2
Set-PctCell $ws "H22" 100
Set-NumCell $ws "J22" 21
Set-PctCell $ws "K22" -23.809523809523

# ---- Row 23 (Housing) ----
Set-NumCell $ws "D23" 1
Set-PctCell $ws "E23" 0
Set-NumCell $ws "J23" 34
Set-PctCell $ws "K23" -17.647058823529
Set-PctCell $ws "M23" 16.666666666666

# ---- Row 24 (Petit Larceny) ----
Set-NumCell $ws "C24" 17
Set-NumCell $ws "D24" 21
Set-PctCell $ws "E24" -19.047619047619
Set-NumCell $ws "F24" 74
Set-NumCell $ws "G24" 99
Set-PctCell $ws "H24" -25.252525252525
Set-NumCell $ws "I24" 1023
Set-NumCell $ws "J24" 941
Set-PctCell $ws "K24" 8.714133900106
Set-PctCell $ws "L24" -10.420315236427
Set-PctCell $ws "M24" 38.617886178861

# ---- Row 25 (Misd. Assault) ----
Set-NumCell $ws "C25" 8
Set-NumCell $ws "D25" 9
Set-PctCell $ws "E25" -11.111111111111
Set-NumCell $ws "F25" 38
Set-NumCell $ws "G25" 46
Set-PctCell $ws "H25" -17.391304347826
Set-NumCell $ws "I25" 568
Set-NumCell $ws "J25" 512
Set-PctCell $ws "K25" 10.9375
Set-PctCell $ws "L25" 20.338983050847
Set-PctCell $ws "M25" -24.867724867724

# ---- Row 26 (UCR Rape*) ----
Set-TextZero $ws "F26"
Set-NumCell $ws "G26" 4
Set-PctCell $ws "H26" -100

# ---- Row 27 (Other Sex Crimes) ----
Set-TextZero $ws "C27"
Set-NumCell $ws "D27" 2
Set-PctCell $ws "E27" -100
Set-NumCell $ws "F27" 1
Set-PctCell $ws "H27" -80
Set-NumCell $ws "J27" 51
Set-PctCell $ws "K27" 23.529411764705
Set-PctCell $ws "L27" 5

# ---- Row 28 (Shooting Vic.) ----
Set-TextZero $ws "C28"
Set-TextZero $ws "D28"
Set-TextStar $ws "E28"
Set-PctCell $ws "N28" -92.477876106194

# ---- Row 29 (Shooting Inc.) ----
Set-TextZero $ws "C29"
Set-TextZero $ws "D29"
Set-TextStar $ws "E29"
Set-PctCell $ws "N29" -92.344497607655
